$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "班尼特"
$ws.Range("B3").Value = "班尼特"
$ws.Range("C3").Value = "90+"
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = "9,9,10"
$ws.Range("F3").Value = "天空之刃"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "宗室"
$ws.Range("I3").Value = "宗室"

for ($col = 10; $col -le 42; $col++) {
    $ws.Cells.Item(3, $col).Value = 0
}

$ws.Range("F7").Select()
